$d = $word.ActiveDocument

# Target color 00B050 (RGB) packed as a Word "long" color value (R + G*256 + B*65536)
$newColor = 5287936

function Set-ParagraphColorByText($doc, $searchText) {
    $rng = $doc.Content.Duplicate
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, `
                                $true, 1, $false, "", 0)
    if ($found) {
        $matchStart = $rng.Start
        $matchEnd = $rng.End
        foreach ($p in $doc.Paragraphs) {
            if ($p.Range.Start -le $matchStart -and $p.Range.End -ge $matchEnd) {
                # Set color on the whole paragraph range (includes the paragraph
                # mark) so both the run text and the paragraph-mark run
                # properties (w:pPr/w:rPr) pick up the new color.
                $p.Range.Font.Color = $newColor
                break
            }
        }
    }
}

Set-ParagraphColorByText $d "Create box plots and histograms to understand the distributions and outliers. Perform outlier treatment."
Set-ParagraphColorByText $d "Does the US fare significantly better than the rest of the world in terms of total purchases?"
